$d = $word.ActiveDocument

# 1) "Hello world!" -> "Hello world again!"
$d.Content.Find.Execute("Hello world!", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Hello world again!", 2)

# 2) Append two new paragraphs at the end of the document:
#    - "This is a second paragraph." + " This text is being added to the second
#      paragraph." as two separate runs in one paragraph
#    - "This is a yet another paragraph." as its own paragraph
# Using Range.InsertXML (flat-OPC WordprocessingML) lets us add the two runs
# as genuinely distinct <w:r> elements instead of having them coalesced into
# a single run the way successive Range.InsertAfter calls would.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos - 1, $endPos - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body>' + `
       '<w:p>' + `
       '<w:r><w:t>This is a second paragraph.</w:t></w:r>' + `
       '<w:r><w:t xml:space="preserve"> This text is being added to the second paragraph.</w:t></w:r>' + `
       '</w:p>' + `
       '<w:p><w:r><w:t>This is a yet another paragraph.</w:t></w:r></w:p>' + `
       '</w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml)
